# Refresh crypto symbol price/volume snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells hold text-formatted numbers/percentages (e.g. "308.62", "1.17%"),
# so force text format first to avoid Excel auto-converting these into
# numeric/percentage values.
$cellUpdates = @{
    "D2" = "308.62"
    "E2" = "1.17%"
    "D3" = "36.39"
    "E3" = "1.36%"
    "D4" = "5.054"
    "E4" = "0.81%"
    "D5" = "0.08126"
    "E5" = "0.81%"
    "D6" = "1.993"
    "E6" = "5.76%"
    "E7" = "-0.17%"
    "D8" = "7.862"
    "E8" = "0.26%"
    "D9" = "0.9266"
    "E9" = "-0.63%"
    "D10" = "0.1475"
    "E10" = "11.61%"
    "D11" = "0.1941"
    "E11" = "1.82%"
    "D12" = "0.09099"
    "E12" = "-1.45%"
    "D13" = "0.03530"
    "E13" = "0.57%"
    "D14" = "0.09864"
    "E14" = "-0.23%"
    "D15" = "0.001410"
    "E15" = "-1.47%"
    "D16" = "0.006499"
    "E16" = "1.33%"
    "D17" = "3.846"
    "E17" = "5.11%"
    "E18" = "8.20%"
    "D19" = "0.3450"
    "D20" = "0.1312"
    "E20" = "-0.14%"
    "D21" = "4.799"
    "E21" = "-8.37%"
    "D23" = "0.04363"
    "E23" = "-1.34%"
    "D24" = "0.001233"
    "E24" = "-0.14%"
    "D25" = "0.004165"
    "E25" = "-11.71%"
    "D27" = "0.0001302"
    "E27" = "-0.03%"
    "D39" = "0.02127"
    "E39" = "8.78%"
    "D40" = "0.05120"
    "E40" = "-0.66%"
    "D41" = "0.007470"
    "E41" = "-1.08%"
    "D42" = "0.01007"
    "E42" = "-1.32%"
    "D43" = "0.1370"
    "D44" = "0.002132"
    "E44" = "-1.87%"
    "D45" = "0.009704"
    "E45" = "-10.10%"
    "D46" = "0.00006273"
    "E46" = "-1.08%"
    "E47" = "-0.06%"
    "E48" = "-0.64%"
    "D49" = "0.001602"
    "E49" = "-3.60%"
    "E50" = "-0.06%"
    "E51" = "-0.06%"
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
}
